$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.482.73"
$ws.Range("D2").Style = $s
$s = $ws.Range("E2").Style
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.90%  "
$ws.Range("E2").Style = $s
$s = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.967.57"
$ws.Range("D3").Style = $s
$s = $ws.Range("E3").Style
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -5.18%  "
$ws.Range("E3").Style = $s
$s = $ws.Range("E4").Style
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E4").Style = $s
$s = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.63"
$ws.Range("D5").Style = $s
$s = $ws.Range("E5").Style
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.38%  "
$ws.Range("E5").Style = $s
$s = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.612"
$ws.Range("D6").Style = $s
$s = $ws.Range("E6").Style
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.27%  "
$ws.Range("E6").Style = $s
$s = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.47"
$ws.Range("D7").Style = $s
$s = $ws.Range("E7").Style
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -10.59%  "
$ws.Range("E7").Style = $s
$s = $ws.Range("E8").Style
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E8").Style = $s
$s = $ws.Range("E9").Style
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -6.74%  "
$ws.Range("E9").Style = $s
$s = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "55.78"
$ws.Range("D10").Style = $s
$s = $ws.Range("E10").Style
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.54%  "
$ws.Range("E10").Style = $s
$s = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0839"
$ws.Range("D11").Style = $s
$s = $ws.Range("E11").Style
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("E11").Style = $s
$s = $ws.Range("E12").Style
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("E12").Style = $s
$s = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.20"
$ws.Range("D13").Style = $s
$s = $ws.Range("E13").Style
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -6.69%  "
$ws.Range("E13").Style = $s
$s = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.834"
$ws.Range("D14").Style = $s
$s = $ws.Range("E14").Style
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -10.13%  "
$ws.Range("E14").Style = $s
$s = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.252.99"
$ws.Range("D15").Style = $s
$s = $ws.Range("E15").Style
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -5.19%  "
$ws.Range("E15").Style = $s
$s = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.56"
$ws.Range("D16").Style = $s
$s = $ws.Range("E16").Style
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -8.90%  "
$ws.Range("E16").Style = $s
$s = $ws.Range("E17").Style
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -5.74%  "
$ws.Range("E17").Style = $s
$s = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.964.55"
$ws.Range("D18").Style = $s
$s = $ws.Range("E18").Style
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.14%  "
$ws.Range("E18").Style = $s
$s = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.326.74"
$ws.Range("D19").Style = $s
$s = $ws.Range("E19").Style
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.14%  "
$ws.Range("E19").Style = $s
$s = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.25"
$ws.Range("D20").Style = $s
$s = $ws.Range("E20").Style
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.32%  "
$ws.Range("E20").Style = $s
$s = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0883"
$ws.Range("D21").Style = $s
$s = $ws.Range("E21").Style
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.04%  "
$ws.Range("E21").Style = $s
$s = $ws.Range("E22").Style
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -6.84%  "
$ws.Range("E22").Style = $s
$s = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.73"
$ws.Range("D23").Style = $s
$s = $ws.Range("E23").Style
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.36%  "
$ws.Range("E23").Style = $s
$s = $ws.Range("E24").Style
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E24").Style = $s
$s = $ws.Range("E25").Style
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.86%  "
$ws.Range("E25").Style = $s
$s = $ws.Range("E26").Style
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.13%  "
$ws.Range("E26").Style = $s
$s = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.63"
$ws.Range("D27").Style = $s
$s = $ws.Range("E27").Style
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.36%  "
$ws.Range("E27").Style = $s
$s = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.15"
$ws.Range("D28").Style = $s
$s = $ws.Range("E28").Style
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("E28").Style = $s
$s = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.91"
$ws.Range("D29").Style = $s
$s = $ws.Range("E29").Style
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.28%  "
$ws.Range("E29").Style = $s
$s = $ws.Range("E30").Style
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.59%  "
$ws.Range("E30").Style = $s
$s = $ws.Range("E31").Style
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.40%  "
$ws.Range("E31").Style = $s
$s = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.16"
$ws.Range("D32").Style = $s
$s = $ws.Range("E32").Style
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.14%  "
$ws.Range("E32").Style = $s
$s = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.76"
$ws.Range("D33").Style = $s
$s = $ws.Range("E33").Style
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -8.53%  "
$ws.Range("E33").Style = $s
$s = $ws.Range("E34").Style
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("E34").Style = $s
$s = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.33"
$ws.Range("D35").Style = $s
$s = $ws.Range("E35").Style
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -8.03%  "
$ws.Range("E35").Style = $s
$s = $ws.Range("E36").Style
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E36").Style = $s
$s = $ws.Range("E37").Style
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("E37").Style = $s
$s = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.98"
$ws.Range("D38").Style = $s
$s = $ws.Range("E38").Style
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -7.48%  "
$ws.Range("E38").Style = $s
$s = $ws.Range("E39").Style
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -14.84%  "
$ws.Range("E39").Style = $s
$s = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("D40").Style = $s
$s = $ws.Range("E40").Style
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.86%  "
$ws.Range("E40").Style = $s
$s = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0966"
$ws.Range("D41").Style = $s
$s = $ws.Range("E41").Style
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -6.64%  "
$ws.Range("E41").Style = $s
$s = $ws.Range("E42").Style
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.69%  "
$ws.Range("E42").Style = $s
$s = $ws.Range("E43").Style
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -8.34%  "
$ws.Range("E43").Style = $s
$s = $ws.Range("E44").Style
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.10%  "
$ws.Range("E44").Style = $s
$s = $ws.Range("E45").Style
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -9.77%  "
$ws.Range("E45").Style = $s
$s = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.79"
$ws.Range("D46").Style = $s
$s = $ws.Range("E46").Style
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -8.80%  "
$ws.Range("E46").Style = $s
$s = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.06"
$ws.Range("D47").Style = $s
$s = $ws.Range("E47").Style
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -7.07%  "
$ws.Range("E47").Style = $s
$s = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.350.13"
$ws.Range("D48").Style = $s
$s = $ws.Range("E48").Style
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.68%  "
$ws.Range("E48").Style = $s
$s = $ws.Range("E49").Style
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -9.15%  "
$ws.Range("E49").Style = $s
$s = $ws.Range("E50").Style
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.57%  "
$ws.Range("E50").Style = $s
$s = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.85"
$ws.Range("D51").Style = $s
$s = $ws.Range("E51").Style
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.18%  "
$ws.Range("E51").Style = $s
